# edit.ps1
#
# Applies the "Add MODE OF OPERATION table" change:
#   1. In the existing OP-CODE table, the "MUL" cell gains a second run
#      containing "T" so the visible text becomes "MULT".
#   2. A new 3x11 "Mode / Code / Description" table (plus a trailing
#      blank paragraph) is inserted right after the bold
#      "MODE OF OPERATION" heading paragraph.

$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part 1: "MUL" -> "MUL" + "T" (two runs, same rPr) in the OP-CODE table.
# -----------------------------------------------------------------------

$mulRange = $d.Content
$mulFound = $mulRange.Find.Execute("MUL", $true, $true, $false, $false, $false, `
                                    $true, 1, $false, "", 0)

if ($mulFound) {
    # Collapse to the end of "MUL" and append "T". Word's model merges
    # the appended text into the same run as the preceding text when the
    # formatting is identical, giving a single run reading "MULT".
    $mulRange.Collapse(0)
    $mulRange.InsertAfter("T")

    # Find the combined "MULT" text again and isolate its final
    # character. Toggling that single character's font size (away from,
    # then back to, the original 12pt/sz=24) forces the run to split
    # into two runs without changing the rendered formatting, matching
    # the diff's extra <w:r> that wraps just the "T".
    $multRange = $d.Content
    $multFound = $multRange.Find.Execute("MULT", $true, $true, $false, $false, $false, `
                                          $true, 1, $false, "", 0)
    if ($multFound) {
        $lastChar = $d.Range($multRange.End - 1, $multRange.End)
        $lastChar.Font.Size = 13
        $lastChar.Font.Size = 12
    }
}

# -----------------------------------------------------------------------
# Part 2: build the new "MODE OF OPERATION" table as a WordOpenXML blob.
# -----------------------------------------------------------------------

function Get-CellXml([string]$text, [int]$width, [bool]$bold) {
    if ($bold) {
        $rPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
    } else {
        $rPr = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
    }
    return "<w:tc><w:tcPr><w:tcW w:w=""$width"" w:type=""dxa""/></w:tcPr>" +
           "<w:p><w:pPr>$rPr</w:pPr><w:r>$rPr<w:t>$text</w:t></w:r></w:p></w:tc>"
}

function Get-RowXml($cols, [bool]$bold, $widths) {
    $tr = "<w:tr>"
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $tr += Get-CellXml $cols[$i] $widths[$i] $bold
    }
    $tr += "</w:tr>"
    return $tr
}

$colWidths = @(3116, 3117, 3117)

$tableRows = @(
    @("Mode", "Code", "Description"),
    @("Ready", "0000", "System is Idle"),
    @("ADD", "0001", "Perform Addition"),
    @("SUB", "0010", "Perform Subtraction"),
    @("MULT", "0011", "Perform Multiplication"),
    @("DIV", "0100", "Perform Division"),
    @("AND", "0101", "Perform AND operation"),
    @("OR", "0110", "Perform OR operation"),
    @("XOR", "0111", "Perform XOR operation"),
    @("NOT", "1000", "Perform NOT operation"),
    @("RESET", "1111", "Reset ACC with 0s")
)

$rowsXml = ""
for ($r = 0; $r -lt $tableRows.Length; $r++) {
    $isHeader = ($r -eq 0)
    $rowsXml += Get-RowXml $tableRows[$r] $isHeader $colWidths
}

$tblPr = '<w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/>' +
         '<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" ' +
         'w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr>'
$tblGrid = '<w:tblGrid><w:gridCol w:w="3116"/><w:gridCol w:w="3117"/><w:gridCol w:w="3117"/></w:tblGrid>'

$tableXmlFragment = "<w:tbl>$tblPr$tblGrid$rowsXml</w:tbl>" +
                     '<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'

$wordOpenXml = '<?xml version="1.0" standalone="yes"?>' +
    '<?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" ' +
    'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    "<w:body>$tableXmlFragment</w:body>" +
    '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# -----------------------------------------------------------------------
# Part 3: locate the "MODE OF OPERATION" heading paragraph and insert the
# table (+ trailing blank paragraph) immediately after it.
# -----------------------------------------------------------------------

$headingRange = $d.Content
$headingFound = $headingRange.Find.Execute("MODE OF OPERATION", $true, $true, $false, `
                                            $false, $false, $true, 1, $false, "", 0)

if ($headingFound) {
    $headingPara = $headingRange.Paragraphs.First

    # Using a zero-length Range positioned right at the paragraph's end
    # (rather than Collapse()-ing the Find range itself) keeps the
    # "MODE OF OPERATION" text intact -- InsertXML on a collapsed Find
    # range can otherwise swallow the preceding paragraph.
    $insertPoint = $d.Range($headingPara.Range.End, $headingPara.Range.End)
    [void]$insertPoint.InsertXML($wordOpenXml)
}
